# Natmi following Dr Hou advice
#
# The natmi ligand-receptor (Tnf -> Tnfrsf1b) edge table was regenerated
# with an additional "FAPs" sending cluster. Existing ECs -> {ECs,FAPs,sCs}
# rows get refreshed statistics (row 2-4), and three new rows are appended
# for FAPs -> {ECs,FAPs,sCs} (row 5-7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

  # Row 2: ECs -> Tnf/Tnfrsf1b -> ECs
  $ws.Range("A2").Value = "ECs"
  $ws.Range("B2").Value = "Tnf"
  $ws.Range("C2").Value = "Tnfrsf1b"
  $ws.Range("D2").Value = "ECs"
  $ws.Range("E2").Value = 3
  $ws.Range("F2").Value = 1
  $ws.Range("G2").Value = 459.5553626666667
  $ws.Range("H2").Value = 1378.666088
  $ws.Range("I2").Value = 0.997342491800172
  $ws.Range("J2").Value = 0.997342491800172
  $ws.Range("K2").Value = 3
  $ws.Range("L2").Value = 1
  $ws.Range("M2").Value = 87.481206
  $ws.Range("N2").Value = 262.443618
  $ws.Range("O2").Value = 0.8890954013316028
  $ws.Range("P2").Value = 0.8890954013316029
  $ws.Range("Q2").Value = 40202.45734984738
  $ws.Range("R2").Value = 361822.1161486264
  $ws.Range("S2").Value = 0.8867326230121347
  $ws.Range("T2").Value = 0.8867326230121348
  # Row 3: ECs -> Tnf/Tnfrsf1b -> FAPs
  $ws.Range("A3").Value = "ECs"
  $ws.Range("B3").Value = "Tnf"
  $ws.Range("C3").Value = "Tnfrsf1b"
  $ws.Range("D3").Value = "FAPs"
  $ws.Range("E3").Value = 3
  $ws.Range("F3").Value = 1
  $ws.Range("G3").Value = 459.5553626666667
  $ws.Range("H3").Value = 1378.666088
  $ws.Range("I3").Value = 0.997342491800172
  $ws.Range("J3").Value = 0.997342491800172
  $ws.Range("K3").Value = 3
  $ws.Range("L3").Value = 1
  $ws.Range("M3").Value = 10.36557633333333
  $ws.Range("N3").Value = 31.096729
  $ws.Range("O3").Value = 0.1053481847303107
  $ws.Range("P3").Value = 0.1053481847303107
  $ws.Range("Q3").Value = 4763.556191114017
  $ws.Range("R3").Value = 42872.00572002615
  $ws.Range("S3").Value = 0.1050682210655529
  $ws.Range("T3").Value = 0.1050682210655529
  # Row 4: ECs -> Tnf/Tnfrsf1b -> sCs
  $ws.Range("A4").Value = "ECs"
  $ws.Range("B4").Value = "Tnf"
  $ws.Range("C4").Value = "Tnfrsf1b"
  $ws.Range("D4").Value = "sCs"
  $ws.Range("E4").Value = 3
  $ws.Range("F4").Value = 1
  $ws.Range("G4").Value = 459.5553626666667
  $ws.Range("H4").Value = 1378.666088
  $ws.Range("I4").Value = 0.997342491800172
  $ws.Range("J4").Value = 0.997342491800172
  $ws.Range("K4").Value = 3
  $ws.Range("L4").Value = 1
  $ws.Range("M4").Value = 0.546715
  $ws.Range("N4").Value = 1.640145
  $ws.Range("O4").Value = 0.005556413938086396
  $ws.Range("P4").Value = 0.005556413938086396
  $ws.Range("Q4").Value = 251.2458101003066
  $ws.Range("R4").Value = 2261.21229090276
  $ws.Range("S4").Value = 0.005541647722484292
  $ws.Range("T4").Value = 0.005541647722484292
  # Row 5: FAPs -> Tnf/Tnfrsf1b -> ECs
  $ws.Range("A5").Value = "FAPs"
  $ws.Range("B5").Value = "Tnf"
  $ws.Range("C5").Value = "Tnfrsf1b"
  $ws.Range("D5").Value = "ECs"
  $ws.Range("E5").Value = 2
  $ws.Range("F5").Value = 0.6666666666666666
  $ws.Range("G5").Value = 1.224526333333333
  $ws.Range("H5").Value = 3.673579
  $ws.Range("I5").Value = 0.002657508199827995
  $ws.Range("J5").Value = 0.002657508199827995
  $ws.Range("K5").Value = 3
  $ws.Range("L5").Value = 1
  $ws.Range("M5").Value = 87.481206
  $ws.Range("N5").Value = 262.443618
  $ws.Range("O5").Value = 0.8890954013316028
  $ws.Range("P5").Value = 0.8890954013316029
  $ws.Range("Q5").Value = 107.123040418758
  $ws.Range("R5").Value = 964.1073637688221
  $ws.Range("S5").Value = 0.002362778319468097
  $ws.Range("T5").Value = 0.002362778319468097
  # Row 6: FAPs -> Tnf/Tnfrsf1b -> FAPs
  $ws.Range("A6").Value = "FAPs"
  $ws.Range("B6").Value = "Tnf"
  $ws.Range("C6").Value = "Tnfrsf1b"
  $ws.Range("D6").Value = "FAPs"
  $ws.Range("E6").Value = 2
  $ws.Range("F6").Value = 0.6666666666666666
  $ws.Range("G6").Value = 1.224526333333333
  $ws.Range("H6").Value = 3.673579
  $ws.Range("I6").Value = 0.002657508199827995
  $ws.Range("J6").Value = 0.002657508199827995
  $ws.Range("K6").Value = 3
  $ws.Range("L6").Value = 1
  $ws.Range("M6").Value = 10.36557633333333
  $ws.Range("N6").Value = 31.096729
  $ws.Range("O6").Value = 0.1053481847303107
  $ws.Range("P6").Value = 0.1053481847303107
  $ws.Range("Q6").Value = 12.69292118034345
  $ws.Range("R6").Value = 114.236290623091
  $ws.Range("S6").Value = 0.000279963664757795
  $ws.Range("T6").Value = 0.0002799636647577951
  # Row 7: FAPs -> Tnf/Tnfrsf1b -> sCs
  $ws.Range("A7").Value = "FAPs"
  $ws.Range("B7").Value = "Tnf"
  $ws.Range("C7").Value = "Tnfrsf1b"
  $ws.Range("D7").Value = "sCs"
  $ws.Range("E7").Value = 2
  $ws.Range("F7").Value = 0.6666666666666666
  $ws.Range("G7").Value = 1.224526333333333
  $ws.Range("H7").Value = 3.673579
  $ws.Range("I7").Value = 0.002657508199827995
  $ws.Range("J7").Value = 0.002657508199827995
  $ws.Range("K7").Value = 3
  $ws.Range("L7").Value = 1
  $ws.Range("M7").Value = 0.546715
  $ws.Range("N7").Value = 1.640145
  $ws.Range("O7").Value = 0.005556413938086396
  $ws.Range("P7").Value = 0.005556413938086396
  $ws.Range("Q7").Value = 0.6694669143283333
  $ws.Range("R7").Value = 6.025202228955
  $ws.Range("S7").Value = 0.00001476621560210316
  $ws.Range("T7").Value = 0.00001476621560210316
